$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round row 5 data values to 2 decimal places (custom accuracy)
$ws.Range("B5").Value = 8.17
$ws.Range("C5").Value = 5.79
$ws.Range("D5").Value = 0.83
$ws.Range("E5").Value = 17.72
$ws.Range("F5").Value = 14.24
$ws.Range("G5").Value = 6.43
$ws.Range("H5").Value = 27.68
$ws.Range("I5").Value = 9.89
$ws.Range("J5").Value = 4.3
$ws.Range("K5").Value = 6.23
$ws.Range("L5").Value = 7.11
$ws.Range("M5").Value = 7.44
$ws.Range("N5").Value = 2.06
$ws.Range("O5").Value = 6.39
$ws.Range("P5").Value = 9.02
$ws.Range("Q5").Value = 5.56
$ws.Range("R5").Value = 0.74
$ws.Range("S5").Value = 0.48
$ws.Range("T5").Value = 89.97
$ws.Range("U5").Value = 17.99
$ws.Range("V5").Value = 5.9
$ws.Range("W5").Value = 11.91
$ws.Range("X5").Value = 6.18
$ws.Range("Y5").Value = 1.16
$ws.Range("Z5").Value = 13.05
$ws.Range("AA5").Value = 5.21
$ws.Range("AB5").Value = 4.73
$ws.Range("AC5").Value = 5.54
$ws.Range("AD5").Value = 7.41
$ws.Range("AE5").Value = 0.56
$ws.Range("AF5").Value = 25.26
$ws.Range("AG5").Value = 3.23
$ws.Range("AH5").Value = 7.38

# Remove the last data row (row 6) - reducing dataset, trimming to 1000 rows worth of data
$ws.Rows.Item(6).Delete()
